# Append the 09/09/2025 profit allocation row (row 8) to Sheet1,
# matching the existing table of Date / BTC / KAS fraction columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the date as plain text (e.g. "09/08/2025" on the row
# above), not an actual date value, so force a text format before
# assigning the value to stop Excel from auto-converting it to a date
# serial number. ClearFormats() afterwards drops the temporary "@"
# number-format style again so the new cell doesn't end up with an
# explicit style index that the other data rows don't have.
$dateCell = $ws.Range("A8")
$dateCell.NumberFormat = "@"
$dateCell.Value = "09/09/2025"
$dateCell.ClearFormats()

$ws.Range("B8").Value = 0.1246141075485167
$ws.Range("C8").Value = 0.8753858924514833
